$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value and force it to remain plain text even if it
# looks like a number (Excel would otherwise auto-convert "600.93" etc.
# into a numeric value). We briefly mark the cell as Text, assign the
# value, then restore the default "Normal" style so no visible/style
# change is left behind on the cell.
function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '64.664.77'
$ws.Cells.Item(2, 5).Value = '  -2.99%  '
$ws.Cells.Item(3, 4).Value = '3.166.13'
$ws.Cells.Item(3, 5).Value = '  -1.98%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
Set-TextValue 5 4 '600.93'
$ws.Cells.Item(5, 5).Value = '  -1.41%  '
Set-TextValue 6 4 '151.02'
$ws.Cells.Item(6, 5).Value = '  -5.08%  '
$ws.Cells.Item(7, 5).Value = '  +0.06%  '
$ws.Cells.Item(8, 4).Value = '3.166.61'
$ws.Cells.Item(8, 5).Value = '  -1.96%  '
Set-TextValue 9 4 '0.536'
$ws.Cells.Item(9, 5).Value = '  -2.89%  '
$ws.Cells.Item(10, 5).Value = '  -4.37%  '
Set-TextValue 11 4 '5.66'
$ws.Cells.Item(11, 5).Value = '  -1.19%  '
$ws.Cells.Item(12, 5).Value = '  -4.74%  '
$ws.Cells.Item(13, 5).Value = '  -3.91%  '
Set-TextValue 14 4 '37.21'
$ws.Cells.Item(14, 5).Value = '  -4.26%  '
$ws.Cells.Item(15, 4).Value = '3.664.90'
$ws.Cells.Item(15, 5).Value = '  -2.53%  '
$ws.Cells.Item(16, 4).Value = '64.731.35'
$ws.Cells.Item(16, 5).Value = '  -2.97%  '
$ws.Cells.Item(17, 5).Value = '  +0.48%  '
$ws.Cells.Item(18, 4).Value = '3.166.30'
$ws.Cells.Item(18, 5).Value = '  -2.15%  '
$ws.Cells.Item(19, 5).Value = '  -4.21%  '
Set-TextValue 20 4 '484.35'
$ws.Cells.Item(20, 5).Value = '  -4.77%  '
Set-TextValue 21 4 '14.89'
$ws.Cells.Item(21, 5).Value = '  -2.14%  '
$ws.Cells.Item(22, 5).Value = '  -2.10%  '
$ws.Cells.Item(23, 5).Value = '  -2.39%  '
Set-TextValue 24 4 '14.04'
$ws.Cells.Item(24, 5).Value = '  -3.64%  '
Set-TextValue 25 4 '85.42'
$ws.Cells.Item(25, 5).Value = '  +0.56%  '
$ws.Cells.Item(26, 5).Value = '  -0.06%  '
$ws.Cells.Item(27, 5).Value = '  -2.35%  '
$ws.Cells.Item(28, 5).Value = '  -3.76%  '
$ws.Cells.Item(29, 5).Value = '  -3.91%  '
Set-TextValue 30 4 '7.18'
$ws.Cells.Item(30, 5).Value = '  +2.46%  '
$ws.Cells.Item(31, 5).Value = '  -0.48%  '
$ws.Cells.Item(32, 5).Value = '  -6.43%  '
$ws.Cells.Item(33, 5).Value = '  -0.17%  '
Set-TextValue 34 4 '26.96'
$ws.Cells.Item(34, 5).Value = '  -4.35%  '
$ws.Cells.Item(35, 5).Value = '  -5.26%  '
$ws.Cells.Item(36, 5).Value = '  -4.87%  '
Set-TextValue 37 4 '54.94'
$ws.Cells.Item(37, 5).Value = '  -1.39%  '
Set-TextValue 38 4 '3.26'
$ws.Cells.Item(38, 5).Value = '  +5.64%  '
$ws.Cells.Item(39, 4).Value = '0.0₃0751'
$ws.Cells.Item(39, 5).Value = '  -2.67%  '
Set-TextValue 40 4 '462.96'
$ws.Cells.Item(40, 5).Value = '  -7.56%  '
$ws.Cells.Item(41, 5).Value = '  -3.49%  '
$ws.Cells.Item(42, 5).Value = '  -3.73%  '
Set-TextValue 43 4 '8.58'
$ws.Cells.Item(43, 5).Value = '  -1.53%  '
$ws.Cells.Item(44, 5).Value = '  -0.48%  '
$ws.Cells.Item(45, 4).Value = '2.900.37'
$ws.Cells.Item(45, 5).Value = '  -0.17%  '
$ws.Cells.Item(46, 5).Value = '  -6.81%  '
Set-TextValue 47 4 '27.15'
$ws.Cells.Item(47, 5).Value = '  -3.51%  '
$ws.Cells.Item(48, 5).Value = '  -0.04%  '
$ws.Cells.Item(49, 5).Value = '  -2.77%  '
$ws.Cells.Item(50, 5).Value = '  +0.26%  '

# Row 51: coin changed from Arweave to Monero
$ws.Cells.Item(51, 2).Value = 'Monero'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 51 4 '119.82'
$ws.Cells.Item(51, 5).Value = '  -2.07%  '

